# "Generate Report for Archive"
#
# The localization status report moved from "Ready for handoff" to
# "In Translation" for the two rows tracked in this workbook. That status
# string shows up on the "Overview" sheet (once per locale column, zh-cn
# and de-de) and on each per-locale detail sheet's "Status" column.
#
# Shrinking the status text also lets the Status column(s) narrow, so we
# resize those columns to match.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: per-locale status columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# --- "zh-cn" sheet: Status column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- "de-de" sheet: Status column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C1").ColumnWidth = 12.5
